# Add a new "jour 4" column (F) to the task-repartition sheet, with the same
# header/body styling as the existing columns, then update the view/selection
# state to match where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F content -------------------------------------------------
# Cell values are written in the same row order as the new shared strings
# were introduced so they end up appended to sharedStrings.xml in the
# expected order: "Tâches  jour 4", "Sécurité+Front", "websocket+heroku",
# "partie responsable/sessions+front".
$ws.Range("F2").Value = "Tâches  jour 4"
$ws.Range("F5").Value = "Sécurité+Front"
$ws.Range("F3").Value = "websocket+heroku"
$ws.Range("F4").Value = "partie responsable/sessions+front"

# --- Formatting -------------------------------------------------------------
# Copy the formatting (borders etc.) from column E onto the new column F so
# the look matches the rest of the table: row 2 uses the thick-bordered
# header style, row 3 uses the "top-open" body style, rows 4-5 use the full
# thin-bordered body style.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null

$ws.Range("E3").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null

$ws.Range("E4").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null

$ws.Range("E5").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Column width ------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 35.92

# --- View state --------------------------------------------------------------
$ws.Range("F4").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
